# NYPD 45th Precinct CompStat weekly update
# - Bumps the report volume/week number and the covered date range
# - Refreshes this week's crime-complaint figures (rows 14-29)
#
# Some cells flip between a numeric value and a literal text placeholder
# ("0" / "***.*") that the report uses when a percentage change is not
# meaningful (e.g. division by zero). Those placeholders are plain text
# cells (not real blanks/zeros) in the sheet, with their own cell style.
# We use Range.Copy() from a cell that already carries the desired
# style+content to flip a cell's type/style cleanly, then (for numeric
# targets) overwrite the copied value with the correct number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: volume/number and the covered week's date range
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# ---------------------------------------------------------------------
# Helper donor cells that keep a stable style/content for the whole
# script (never themselves become the "odd one out" placeholder type):
#   C14 -> style 14, text "0"     (shared string 20)
#   M14 -> style 14, text "***.*" (shared string 21)
#   I22 -> style 15 (plain count number format)
#   K22 -> style 16 (percent-change number format)
# ---------------------------------------------------------------------

# Row 14 (Murder): D14, E14 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("M14").Copy($ws.Range("E14"))

# Row 15 (Rape): C15 switches from number to text placeholder "0"
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -36.842105263157

# Row 16 (Robbery)
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("I16").Value = 148
$ws.Range("J16").Value = 145
$ws.Range("K16").Value = 2.068965517241
$ws.Range("L16").Value = 13.846153846153
$ws.Range("M16").Value = -7.5
$ws.Range("N16").Value = -58.192090395480

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 222
$ws.Range("J17").Value = 197
$ws.Range("K17").Value = 12.690355329949
$ws.Range("L17").Value = 35.365853658536
$ws.Range("M17").Value = 56.338028169014
$ws.Range("N17").Value = 26.136363636363

# Row 18 (Burglary)
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = 34.210526315789
$ws.Range("L18").Value = 29.113924050632
$ws.Range("M18").Value = -56.595744680851
$ws.Range("N18").Value = -84.186046511627

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 36
$ws.Range("H19").Value = -32.075471698113
$ws.Range("I19").Value = 476
$ws.Range("J19").Value = 402
$ws.Range("K19").Value = 18.407960199005
$ws.Range("L19").Value = 43.806646525679
$ws.Range("M19").Value = 37.572254335260
$ws.Range("N19").Value = 49.216300940438

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = -2.857142857142
$ws.Range("I20").Value = 394
$ws.Range("J20").Value = 219
$ws.Range("K20").Value = 79.908675799086
$ws.Range("L20").Value = 84.976525821596
$ws.Range("M20").Value = 146.25
$ws.Range("N20").Value = -74.382314694408

# Row 21 (TOTAL)
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 7.142857142857
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = -14.814814814814
$ws.Range("I21").Value = 1357
$ws.Range("J21").Value = 1055
$ws.Range("K21").Value = 28.625592417061
$ws.Range("L21").Value = 45.757250268528
$ws.Range("M21").Value = 28.139754485363
$ws.Range("N21").Value = -55.682560418027

# Row 22 (Transit): C22,D22,F22,G22,E22,H22 switch from text placeholders to numbers
$ws.Range("I22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("I22").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K22").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 0
$ws.Range("I22").Copy($ws.Range("F22"))
$ws.Range("F22").Value = 1
$ws.Range("I22").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 1
$ws.Range("K22").Copy($ws.Range("H22"))
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = -16.666666666666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0

# Row 23 (Housing)
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("I23").Value = 47
$ws.Range("K23").Value = 88
$ws.Range("L23").Value = 38.235294117647
$ws.Range("M23").Value = 30.555555555555

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 21.052631578947
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 21.978021978022
$ws.Range("I24").Value = 1037
$ws.Range("J24").Value = 855
$ws.Range("K24").Value = 21.286549707602
$ws.Range("L24").Value = 50.946142649199
$ws.Range("M24").Value = -5.210237659963

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -3.030303030303
$ws.Range("I25").Value = 372
$ws.Range("J25").Value = 333
$ws.Range("K25").Value = 11.711711711711
$ws.Range("L25").Value = 16.981132075471
$ws.Range("M25").Value = 17.350157728706

# Row 26 (UCR Rape*): C26,D26,E26 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("M14").Copy($ws.Range("E26"))
$ws.Range("I26").Value = 25
$ws.Range("K26").Value = -3.846153846153
$ws.Range("L26").Value = 92.307692307692

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -20.512820512820
$ws.Range("L27").Value = 0

# Row 28 (Shooting Vic.): D28, E28 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("M14").Copy($ws.Range("E28"))

# Row 29 (Shooting Inc.): D29, E29 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("M14").Copy($ws.Range("E29"))
